$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date number format / style) from the cell above onto G8
# before writing values, so it reuses the existing style index instead of
# creating a brand new numFmt entry.
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A8").Value = 9309.17
$ws.Range("B8").Value = 9407
$ws.Range("C8").Value = 109.08
$ws.Range("D8").Value = 107.95
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -1.04
$ws.Range("G8").Value = 42612.672962962963
$ws.Range("H8").Value = $false
